$d = $word.ActiveDocument

# =====================================================================
# Change 1: Title paragraph -- wrap the position-paper name in curly
# quotes:
#   "Motion to Adopt a Position Paper Statement on Mental Health and
#    Suicide in Indigenous Communities in Canada"
# becomes
#   "Motion to Adopt a Position Paper "Mental Health and Suicide in
#    Indigenous Communities in Canada""  (curly quotes)
# =====================================================================
$p1 = $d.Paragraphs.Item(1)
$t1 = $p1.Range.Text
$start1 = $p1.Range.Start
$idxStatement = $t1.IndexOf("Statement on ")
$lenStatement = "Statement on ".Length
$delStart = $start1 + $idxStatement
$delEnd = $delStart + $lenStatement
$rStatement = $d.Range($delStart, $delEnd)
$rStatement.Text = [string][char]0x201C

# Re-fetch paragraph 1 and append the closing curly quote after "Canada"
$p1 = $d.Paragraphs.Item(1)
$endOfP1Text = $p1.Range.End - 1
$closeQuotePos = $endOfP1Text
$rClose = $d.Range($closeQuotePos, $closeQuotePos)
$rClose.InsertAfter([string][char]0x201D)
$rCloseFmt = $d.Range($closeQuotePos, $closeQuotePos + 1)
$rCloseFmt.Font.Name = "Times New Roman"
$rCloseFmt.Font.NameFarEast = "Times New Roman"
$rCloseFmt.Font.NameBi = "Times New Roman"
$rCloseFmt.Bold = 1
$rCloseFmt.Font.Size = 13
$rCloseFmt.Font.Underline = 1

# =====================================================================
# Change 2: BE IT RESOLVED paragraph -- insert "Mental Health and "
# before "Suicide in Indigenous Communities in Canada" so the position
# paper title reads in full, keeping the straight double quotes intact.
# =====================================================================
$p9 = $d.Paragraphs.Item(9)
$t9 = $p9.Range.Text
$start9 = $p9.Range.Start
$idxSuicide = $t9.IndexOf("Suicide in Indigenous")
$insertPos9 = $start9 + $idxSuicide
$rSuicide = $d.Range($insertPos9, $insertPos9)
$rSuicide.InsertBefore("Mental Health and ")

# =====================================================================
# Change 3 / 4: Relocate the _GoBack bookmark from the very end of the
# document to immediately after " none" (end of the "Anticipated
# Financial Cost: none" paragraph).
# =====================================================================
$oldGoBack = $d.Bookmarks.Item("_GoBack")
$oldGoBack.Delete()

$pCost = $d.Paragraphs.Item(11)
$goBackPos = $pCost.Range.End - 1
$rTmp = $d.Range($goBackPos, $goBackPos)
$rTmp.InsertAfter("Q")
$rTmpRange = $d.Range($goBackPos, $goBackPos + 1)
$d.Bookmarks.Add("_GoBack", $rTmpRange)
$rTmpDelete = $d.Range($goBackPos, $goBackPos + 1)
$rTmpDelete.Delete()

# =====================================================================
# Change 5: Split the "Moved by: Kai Homer, University of Alberta"
# paragraph into three paragraphs:
#   "Moved by: "
#   "Kai Homer, University of Alberta "
#   "Ali Sumner, University of Toronto"
# =====================================================================
$pMoved = $d.Paragraphs.Item(15)
$tMoved = $pMoved.Range.Text
$startMoved = $pMoved.Range.Start
$idxKai = $tMoved.IndexOf("Kai Homer")
$splitPosKai = $startMoved + $idxKai
$rSplitKai = $d.Range($splitPosKai, $splitPosKai)
$rSplitKai.InsertParagraphAfter()

# Paragraph with "Kai Homer, University of Alberta" -- add trailing space
$pKai = $d.Paragraphs.Item(16)
$endKai = $pKai.Range.End - 1
$rSpace = $d.Range($endKai, $endKai)
$rSpace.InsertAfter(" ")
$rSpaceFmt = $d.Range($endKai, $endKai + 1)
$rSpaceFmt.Font.Name = "Times New Roman"
$rSpaceFmt.Font.Size = 14

# Split again after the trailing space to create a new, empty paragraph
# for "Ali Sumner, University of Toronto"
$pKai = $d.Paragraphs.Item(16)
$endKai2 = $pKai.Range.End - 1
$rSplitAli = $d.Range($endKai2, $endKai2)
$rSplitAli.InsertParagraphAfter()

$pAli = $d.Paragraphs.Item(17)
$aliStart = $pAli.Range.Start
$rAli = $d.Range($aliStart, $aliStart)
$rAli.InsertAfter("Ali Sumner, University of Toronto")
